$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet originally listed 19 exam-takers (rows 2-20, TaiKhoanId
# 1..19) all tied to exam period 2 (column B). The refreshed admin export
# only has 10 exam-takers tied to exam period 1, and TaiKhoanId now starts
# at 2 - i.e. the former row 2 (TaiKhoanId 1) dropped off and everything
# else shifted up one. Delete the old first data row, then trim the tail
# back down to 10 rows (this also shrinks the sheet dimension from
# A1:C20 down to A1:C11). Deleting rows (rather than overwriting the
# values) keeps column A's existing "quote prefix" cell style intact.
$ws.Range("A2:C2").EntireRow.Delete()
$ws.Range("A12:C19").EntireRow.Delete()

# Every remaining exam-taker now belongs to exam period 1 instead of 2.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = 1
}

# Column B (KyKiemTraId) is slightly wider to fit the refreshed layout.
$ws.Columns.Item(2).ColumnWidth = 13.33

# Move the active selection like the resaved workbook shows.
$ws.Range("F7").Select()
